$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 1233.1666  # H4: 1488.3334 -> 1233.1666
$ws.Cells.Item(4, 9).Value = 349.75  # I4: 233.75 -> 349.75
$ws.Cells.Item(4, 10).Value = 3000  # J4: 3997.5 -> 3000
$ws.Cells.Item(4, 11).Value = 349.75  # K4: 233.75 -> 349.75
$ws.Cells.Item(4, 12).Value = 3000  # L4: 3997.5 -> 3000
$ws.Cells.Item(4, 13).Value = -235.75  # M4: -119.75 -> -235.75
$ws.Cells.Item(4, 14).Value = -3228  # N4: -4225.5 -> -3228
$ws.Cells.Item(8, 8).Value = 90  # H8: 50 -> 90
$ws.Cells.Item(8, 9).Value = 90  # I8: 50 -> 90
$ws.Cells.Item(8, 11).Value = 270  # K8: 150 -> 270
$ws.Cells.Item(8, 13).Value = -131  # M8: -11 -> -131
$ws.Cells.Item(17, 8).Value = 443.3784  # H17: 395.58823 -> 443.3784
$ws.Cells.Item(17, 10).Value = 443.3784  # J17: 395.58823 -> 443.3784
$ws.Cells.Item(17, 12).Value = 1330.1352  # L17: 1186.76469 -> 1330.1352
$ws.Cells.Item(17, 14).Value = -1666.1352  # N17: -1522.76469 -> -1666.1352
$ws.Cells.Item(18, 8).Value = 418.83334  # H18: 439.85 -> 418.83334
$ws.Cells.Item(18, 9).Value = 418.83334  # I18: 420.89474 -> 418.83334
$ws.Cells.Item(18, 10).Value = 0  # J18: 800 -> 0
$ws.Cells.Item(18, 11).Value = 418.83334  # K18: 420.89474 -> 418.83334
$ws.Cells.Item(18, 12).Value = 0  # L18: 800 -> 0
$ws.Cells.Item(18, 13).Value = -134.83334  # M18: -136.89474 -> -134.83334
$ws.Cells.Item(18, 14).ClearContents()  # N18 was -1368
$ws.Cells.Item(25, 8).Value = 10000  # H25: 0 -> 10000
$ws.Cells.Item(25, 9).Value = 10000  # I25: 0 -> 10000
$ws.Cells.Item(25, 11).Value = 30000  # K25: 0 -> 30000
$ws.Cells.Item(25, 13).Value = -29826  # M25: None -> -29826
$ws.Cells.Item(43, 8).Value = 4062.8  # H43: 4680.08 -> 4062.8
$ws.Cells.Item(43, 9).Value = 1870.1428  # I43: 2529.2856 -> 1870.1428
$ws.Cells.Item(43, 10).Value = 4730.1304  # J43: 5516.5 -> 4730.1304
$ws.Cells.Item(43, 11).Value = 1870.1428  # K43: 2529.2856 -> 1870.1428
$ws.Cells.Item(43, 12).Value = 4730.1304  # L43: 5516.5 -> 4730.1304
$ws.Cells.Item(43, 13).Value = -1801.1428  # M43: -2460.2856 -> -1801.1428
$ws.Cells.Item(43, 14).Value = -4868.1304  # N43: -5654.5 -> -4868.1304
$ws.Cells.Item(52, 8).Value = 1443.5  # H52: 1444 -> 1443.5
$ws.Cells.Item(52, 9).Value = 1443.5  # I52: 1444 -> 1443.5
$ws.Cells.Item(52, 11).Value = 4330.5  # K52: 4332 -> 4330.5
$ws.Cells.Item(52, 13).Value = -4170.5  # M52: -4172 -> -4170.5
$ws.Cells.Item(99, 8).Value = 76931816  # H99: 83339816 -> 76931816
$ws.Cells.Item(99, 9).Value = 651.8570999999999  # I99: 600.125 -> 651.8570999999999
$ws.Cells.Item(99, 10).Value = 166684830  # J99: 250018260 -> 166684830
$ws.Cells.Item(99, 11).Value = 1955.5713  # K99: 1800.375 -> 1955.5713
$ws.Cells.Item(99, 12).Value = 500054490  # L99: 750054780 -> 500054490
$ws.Cells.Item(99, 13).Value = -457.5712999999998  # M99: -302.375 -> -457.5712999999998
$ws.Cells.Item(99, 14).Value = -500057486  # N99: -750057776 -> -500057486
$ws.Cells.Item(105, 8).Value = 30000  # H105: 156849.67 -> 30000
$ws.Cells.Item(105, 10).Value = 30000  # J105: 156849.67 -> 30000
$ws.Cells.Item(105, 12).Value = 30000  # L105: 156849.67 -> 30000
$ws.Cells.Item(105, 14).Value = -36988  # N105: -163837.67 -> -36988
$ws.Cells.Item(116, 8).Value = 7792.129  # H116: 9520.591 -> 7792.129
$ws.Cells.Item(116, 9).Value = 7493.222  # I116: 8505.066000000001 -> 7493.222
$ws.Cells.Item(116, 10).Value = 8206  # J116: 11696.714 -> 8206
$ws.Cells.Item(116, 11).Value = 7493.222  # K116: 8505.066000000001 -> 7493.222
$ws.Cells.Item(116, 12).Value = 8206  # L116: 11696.714 -> 8206
$ws.Cells.Item(116, 13).Value = -4051.222  # M116: -5063.066000000001 -> -4051.222
$ws.Cells.Item(116, 14).Value = -15090  # N116: -18580.714 -> -15090
$ws.Cells.Item(127, 8).Value = 4455.0356  # H127: 5055.25 -> 4455.0356
$ws.Cells.Item(127, 9).Value = 5114.364  # I127: 6061.1113 -> 5114.364
$ws.Cells.Item(127, 10).Value = 2037.5  # J127: 2037.6666 -> 2037.5
$ws.Cells.Item(127, 11).Value = 15343.092  # K127: 18183.3339 -> 15343.092
$ws.Cells.Item(127, 12).Value = 6112.5  # L127: 6112.9998 -> 6112.5
$ws.Cells.Item(127, 13).Value = -10383.092  # M127: -13223.3339 -> -10383.092
$ws.Cells.Item(127, 14).Value = -16032.5  # N127: -16032.9998 -> -16032.5
$ws.Cells.Item(129, 8).Value = 33334858  # H129: 47621036 -> 33334858
$ws.Cells.Item(129, 9).Value = 1342  # I129: 2248.5 -> 1342
$ws.Cells.Item(129, 10).Value = 66668376  # J129: 66668550 -> 66668376
$ws.Cells.Item(129, 11).Value = 4026  # K129: 6745.5 -> 4026
$ws.Cells.Item(129, 12).Value = 200005128  # L129: 200005650 -> 200005128
$ws.Cells.Item(129, 13).Value = 974  # M129: -1745.5 -> 974
$ws.Cells.Item(129, 14).Value = -200015128  # N129: -200015650 -> -200015128
$ws.Cells.Item(138, 8).Value = 8692.808000000001  # H138: 8226.212 -> 8692.808000000001
$ws.Cells.Item(138, 10).Value = 8692.639999999999  # J138: 8211.5 -> 8692.639999999999
$ws.Cells.Item(138, 12).Value = 26077.92  # L138: 24634.5 -> 26077.92
$ws.Cells.Item(138, 14).Value = -36357.92  # N138: -34914.5 -> -36357.92
$ws.Cells.Item(141, 8).Value = 1666.6  # H141: 1671.0834 -> 1666.6
$ws.Cells.Item(141, 9).Value = 1454.6364  # I141: 1505.3 -> 1454.6364
$ws.Cells.Item(141, 10).Value = 2249.5  # J141: 2500 -> 2249.5
$ws.Cells.Item(141, 11).Value = 4363.9092  # K141: 4515.9 -> 4363.9092
$ws.Cells.Item(141, 12).Value = 6748.5  # L141: 7500 -> 6748.5
$ws.Cells.Item(141, 13).Value = 816.0907999999999  # M141: 664.1000000000004 -> 816.0907999999999
$ws.Cells.Item(141, 14).Value = -17108.5  # N141: -17860 -> -17108.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2147.5107  # H32: 2242.0444 -> 2147.5107
$ws.Cells.Item(32, 9).Value = 689.7632  # I32: 726.94446 -> 689.7632
$ws.Cells.Item(32, 11).Value = 689.7632  # K32: 726.94446 -> 689.7632
$ws.Cells.Item(32, 13).Value = -402.7632  # M32: -439.94446 -> -402.7632
$ws.Cells.Item(61, 8).Value = 4921  # H61: 5142.4863 -> 4921
$ws.Cells.Item(61, 9).Value = 5041.3335  # I61: 5212.5 -> 5041.3335
$ws.Cells.Item(61, 10).Value = 4424.625  # J61: 4842.4287 -> 4424.625
$ws.Cells.Item(61, 11).Value = 5041.3335  # K61: 5212.5 -> 5041.3335
$ws.Cells.Item(61, 12).Value = 4424.625  # L61: 4842.4287 -> 4424.625
$ws.Cells.Item(61, 13).Value = -4829.3335  # M61: -5000.5 -> -4829.3335
$ws.Cells.Item(61, 14).Value = -4848.625  # N61: -5266.4287 -> -4848.625
$ws.Cells.Item(122, 8).Value = 2881.547  # H122: 2620.1177 -> 2881.547
$ws.Cells.Item(122, 9).Value = 2635.06  # I122: 2347.0208 -> 2635.06
$ws.Cells.Item(122, 11).Value = 7905.18  # K122: 7041.062399999999 -> 7905.18
$ws.Cells.Item(122, 13).Value = -5455.18  # M122: -4591.062399999999 -> -5455.18
$ws.Cells.Item(136, 8).Value = 4921  # H136: 5142.4863 -> 4921
$ws.Cells.Item(136, 9).Value = 5041.3335  # I136: 5212.5 -> 5041.3335
$ws.Cells.Item(136, 10).Value = 4424.625  # J136: 4842.4287 -> 4424.625
$ws.Cells.Item(136, 11).Value = 15124.0005  # K136: 15637.5 -> 15124.0005
$ws.Cells.Item(136, 12).Value = 13273.875  # L136: 14527.2861 -> 13273.875
$ws.Cells.Item(136, 13).Value = -12574.0005  # M136: -13087.5 -> -12574.0005
$ws.Cells.Item(136, 14).Value = -18373.875  # N136: -19627.2861 -> -18373.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 11264.565  # H80: 11267.608 -> 11264.565
$ws.Cells.Item(80, 9).Value = 2105  # I80: 2115 -> 2105
$ws.Cells.Item(80, 11).Value = 2105  # K80: 2115 -> 2105
$ws.Cells.Item(80, 13).Value = -1107  # M80: -1117 -> -1107
$ws.Cells.Item(83, 8).Value = 11264.565  # H83: 11267.608 -> 11264.565
$ws.Cells.Item(83, 9).Value = 2105  # I83: 2115 -> 2105
$ws.Cells.Item(83, 11).Value = 10525  # K83: 10575 -> 10525
$ws.Cells.Item(83, 13).Value = -5533  # M83: -5583 -> -5533
$ws.Cells.Item(86, 8).Value = 5419.091  # H86: 3680 -> 5419.091
$ws.Cells.Item(86, 9).Value = 2482.5  # I86: 2317.7693 -> 2482.5
$ws.Cells.Item(86, 10).Value = 13250  # J86: 9583 -> 13250
$ws.Cells.Item(86, 11).Value = 2482.5  # K86: 2317.7693 -> 2482.5
$ws.Cells.Item(86, 12).Value = 13250  # L86: 9583 -> 13250
$ws.Cells.Item(86, 13).Value = -1359.5  # M86: -1194.7693 -> -1359.5
$ws.Cells.Item(86, 14).Value = -15496  # N86: -11829 -> -15496
$ws.Cells.Item(89, 8).Value = 5419.091  # H89: 3680 -> 5419.091
$ws.Cells.Item(89, 9).Value = 2482.5  # I89: 2317.7693 -> 2482.5
$ws.Cells.Item(89, 10).Value = 13250  # J89: 9583 -> 13250
$ws.Cells.Item(89, 11).Value = 12412.5  # K89: 11588.8465 -> 12412.5
$ws.Cells.Item(89, 12).Value = 66250  # L89: 47915 -> 66250
$ws.Cells.Item(89, 13).Value = -6796.5  # M89: -5972.8465 -> -6796.5
$ws.Cells.Item(89, 14).Value = -77482  # N89: -59147 -> -77482

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4317.1665  # H31: 4795.304 -> 4317.1665
$ws.Cells.Item(31, 9).Value = 1056.6428  # I31: 1223.8889 -> 1056.6428
$ws.Cells.Item(31, 10).Value = 5030.4062  # J31: 5331.0166 -> 5030.4062
$ws.Cells.Item(31, 11).Value = 1056.6428  # K31: 1223.8889 -> 1056.6428
$ws.Cells.Item(31, 12).Value = 5030.4062  # L31: 5331.0166 -> 5030.4062
$ws.Cells.Item(31, 13).Value = -761.6428000000001  # M31: -928.8888999999999 -> -761.6428000000001
$ws.Cells.Item(31, 14).Value = -5620.4062  # N31: -5921.0166 -> -5620.4062
$ws.Cells.Item(34, 8).Value = 4317.1665  # H34: 4795.304 -> 4317.1665
$ws.Cells.Item(34, 9).Value = 1056.6428  # I34: 1223.8889 -> 1056.6428
$ws.Cells.Item(34, 10).Value = 5030.4062  # J34: 5331.0166 -> 5030.4062
$ws.Cells.Item(34, 11).Value = 1056.6428  # K34: 1223.8889 -> 1056.6428
$ws.Cells.Item(34, 12).Value = 5030.4062  # L34: 5331.0166 -> 5030.4062
$ws.Cells.Item(34, 13).Value = -854.6428000000001  # M34: -1021.8889 -> -854.6428000000001
$ws.Cells.Item(34, 14).Value = -5434.4062  # N34: -5735.0166 -> -5434.4062
$ws.Cells.Item(74, 8).Value = 22221.75  # H74: 0 -> 22221.75
$ws.Cells.Item(74, 9).Value = 22221  # I74: 0 -> 22221
$ws.Cells.Item(74, 10).Value = 22222  # J74: 0 -> 22222
$ws.Cells.Item(74, 11).Value = 22221  # K74: 0 -> 22221
$ws.Cells.Item(74, 12).Value = 22222  # L74: 0 -> 22222
$ws.Cells.Item(74, 13).Value = -21347  # M74: None -> -21347
$ws.Cells.Item(74, 14).Value = -23970  # N74: None -> -23970
$ws.Cells.Item(77, 8).Value = 22221.75  # H77: 0 -> 22221.75
$ws.Cells.Item(77, 9).Value = 22221  # I77: 0 -> 22221
$ws.Cells.Item(77, 10).Value = 22222  # J77: 0 -> 22222
$ws.Cells.Item(77, 11).Value = 66663  # K77: 0 -> 66663
$ws.Cells.Item(77, 12).Value = 66666  # L77: 0 -> 66666
$ws.Cells.Item(77, 13).Value = -62295  # M77: None -> -62295
$ws.Cells.Item(77, 14).Value = -75402  # N77: None -> -75402
$ws.Cells.Item(99, 8).Value = 1754.3182  # H99: 1588 -> 1754.3182
$ws.Cells.Item(99, 9).Value = 1137.9333  # I99: 1132.55 -> 1137.9333
$ws.Cells.Item(99, 10).Value = 3075.1428  # J99: 3106.1667 -> 3075.1428
$ws.Cells.Item(99, 11).Value = 1137.9333  # K99: 1132.55 -> 1137.9333
$ws.Cells.Item(99, 12).Value = 3075.1428  # L99: 3106.1667 -> 3075.1428
$ws.Cells.Item(99, 13).Value = 360.0667000000001  # M99: 365.45 -> 360.0667000000001
$ws.Cells.Item(99, 14).Value = -6071.1428  # N99: -6102.1667 -> -6071.1428
$ws.Cells.Item(105, 8).Value = 4584.857  # H105: 5524.75 -> 4584.857
$ws.Cells.Item(105, 9).Value = 650  # I105: 550 -> 650
$ws.Cells.Item(105, 10).Value = 9831.333000000001  # J105: 10499.5 -> 9831.333000000001
$ws.Cells.Item(105, 11).Value = 650  # K105: 550 -> 650
$ws.Cells.Item(105, 12).Value = 9831.333000000001  # L105: 10499.5 -> 9831.333000000001
$ws.Cells.Item(105, 13).Value = 1097  # M105: 1197 -> 1097
$ws.Cells.Item(105, 14).Value = -13325.333  # N105: -13993.5 -> -13325.333
$ws.Cells.Item(122, 8).Value = 2125.0908  # H122: 3584 -> 2125.0908
$ws.Cells.Item(122, 9).Value = 1488.8125  # I122: 2450.25 -> 1488.8125
$ws.Cells.Item(122, 10).Value = 3821.8333  # J122: 4339.8335 -> 3821.8333
$ws.Cells.Item(122, 11).Value = 4466.4375  # K122: 7350.75 -> 4466.4375
$ws.Cells.Item(122, 12).Value = 11465.4999  # L122: 13019.5005 -> 11465.4999
$ws.Cells.Item(122, 13).Value = -2016.4375  # M122: -4900.75 -> -2016.4375
$ws.Cells.Item(122, 14).Value = -16365.4999  # N122: -17919.5005 -> -16365.4999
$ws.Cells.Item(126, 8).Value = 1754.3182  # H126: 1588 -> 1754.3182
$ws.Cells.Item(126, 9).Value = 1137.9333  # I126: 1132.55 -> 1137.9333
$ws.Cells.Item(126, 10).Value = 3075.1428  # J126: 3106.1667 -> 3075.1428
$ws.Cells.Item(126, 11).Value = 3413.7999  # K126: 3397.65 -> 3413.7999
$ws.Cells.Item(126, 12).Value = 9225.428400000001  # L126: 9318.500100000001 -> 9225.428400000001
$ws.Cells.Item(126, 13).Value = -943.7999  # M126: -927.6499999999996 -> -943.7999
$ws.Cells.Item(126, 14).Value = -14165.4284  # N126: -14258.5001 -> -14165.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 254609.75  # H68: 221647.56 -> 254609.75
$ws.Cells.Item(68, 9).Value = 2082.3333  # I68: 2042 -> 2082.3333
$ws.Cells.Item(68, 10).Value = 362835.78  # J68: 317725 -> 362835.78
$ws.Cells.Item(68, 11).Value = 6246.999899999999  # K68: 6126 -> 6246.999899999999
$ws.Cells.Item(68, 12).Value = 1088507.34  # L68: 953175 -> 1088507.34
$ws.Cells.Item(68, 13).Value = -5435.999899999999  # M68: -5315 -> -5435.999899999999
$ws.Cells.Item(68, 14).Value = -1090129.34  # N68: -954797 -> -1090129.34
$ws.Cells.Item(71, 8).Value = 254609.75  # H71: 221647.56 -> 254609.75
$ws.Cells.Item(71, 9).Value = 2082.3333  # I71: 2042 -> 2082.3333
$ws.Cells.Item(71, 10).Value = 362835.78  # J71: 317725 -> 362835.78
$ws.Cells.Item(71, 11).Value = 18740.9997  # K71: 18378 -> 18740.9997
$ws.Cells.Item(71, 12).Value = 3265522.02  # L71: 2859525 -> 3265522.02
$ws.Cells.Item(71, 13).Value = -14684.9997  # M71: -14322 -> -14684.9997
$ws.Cells.Item(71, 14).Value = -3273634.02  # N71: -2867637 -> -3273634.02
$ws.Cells.Item(107, 8).Value = 7319.524  # H107: 6722.174 -> 7319.524
$ws.Cells.Item(107, 10).Value = 8738.883  # J107: 7866.3687 -> 8738.883
$ws.Cells.Item(107, 12).Value = 26216.649  # L107: 23599.1061 -> 26216.649
$ws.Cells.Item(107, 14).Value = -30056.649  # N107: -27439.1061 -> -30056.649

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6333.273  # H70: 6362.8125 -> 6333.273
$ws.Cells.Item(70, 9).Value = 5281.4287  # I70: 5215.7 -> 5281.4287
$ws.Cells.Item(70, 10).Value = 8174  # J70: 8274.666999999999 -> 8174
$ws.Cells.Item(70, 11).Value = 5281.4287  # K70: 5215.7 -> 5281.4287
$ws.Cells.Item(70, 12).Value = 8174  # L70: 8274.666999999999 -> 8174
$ws.Cells.Item(70, 13).Value = -5011.4287  # M70: -4945.7 -> -5011.4287
$ws.Cells.Item(70, 14).Value = -8714  # N70: -8814.666999999999 -> -8714
$ws.Cells.Item(73, 8).Value = 6333.273  # H73: 6362.8125 -> 6333.273
$ws.Cells.Item(73, 9).Value = 5281.4287  # I73: 5215.7 -> 5281.4287
$ws.Cells.Item(73, 10).Value = 8174  # J73: 8274.666999999999 -> 8174
$ws.Cells.Item(73, 11).Value = 5281.4287  # K73: 5215.7 -> 5281.4287
$ws.Cells.Item(73, 12).Value = 8174  # L73: 8274.666999999999 -> 8174
$ws.Cells.Item(73, 13).Value = -4345.4287  # M73: -4279.7 -> -4345.4287
$ws.Cells.Item(73, 14).Value = -10046  # N73: -10146.667 -> -10046
$ws.Cells.Item(102, 8).Value = 3417.0676  # H102: 3908.2122 -> 3417.0676
$ws.Cells.Item(102, 9).Value = 2018.228  # I102: 2451.3877 -> 2018.228
$ws.Cells.Item(102, 11).Value = 2018.228  # K102: 2451.3877 -> 2018.228
$ws.Cells.Item(102, 13).Value = -396.2280000000001  # M102: -829.3877000000002 -> -396.2280000000001
$ws.Cells.Item(113, 8).Value = 8947.333000000001  # H113: 8956 -> 8947.333000000001
$ws.Cells.Item(113, 9).Value = 4202.1665  # I113: 4223.8335 -> 4202.1665
$ws.Cells.Item(113, 11).Value = 4202.1665  # K113: 4223.8335 -> 4202.1665
$ws.Cells.Item(113, 13).Value = -2032.1665  # M113: -2053.8335 -> -2032.1665
$ws.Cells.Item(122, 8).Value = 41318.07  # H122: 6101.364 -> 41318.07
$ws.Cells.Item(122, 9).Value = 42088.355  # I122: 5674.8438 -> 42088.355
$ws.Cells.Item(122, 11).Value = 126265.065  # K122: 17024.5314 -> 126265.065
$ws.Cells.Item(122, 13).Value = -123815.065  # M122: -14574.5314 -> -123815.065

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 15122.656  # H61: 16031.094 -> 15122.656
$ws.Cells.Item(61, 9).Value = 13564.846  # I61: 14518.24 -> 13564.846
$ws.Cells.Item(61, 10).Value = 21873.166  # J61: 21434.143 -> 21873.166
$ws.Cells.Item(61, 11).Value = 13564.846  # K61: 14518.24 -> 13564.846
$ws.Cells.Item(61, 12).Value = 21873.166  # L61: 21434.143 -> 21873.166
$ws.Cells.Item(61, 13).Value = -13362.846  # M61: -14316.24 -> -13362.846
$ws.Cells.Item(61, 14).Value = -22277.166  # N61: -21838.143 -> -22277.166
$ws.Cells.Item(113, 8).Value = 15122.656  # H113: 16031.094 -> 15122.656
$ws.Cells.Item(113, 9).Value = 13564.846  # I113: 14518.24 -> 13564.846
$ws.Cells.Item(113, 10).Value = 21873.166  # J113: 21434.143 -> 21873.166
$ws.Cells.Item(113, 11).Value = 13564.846  # K113: 14518.24 -> 13564.846
$ws.Cells.Item(113, 12).Value = 21873.166  # L113: 21434.143 -> 21873.166
$ws.Cells.Item(113, 13).Value = -11394.846  # M113: -12348.24 -> -11394.846
$ws.Cells.Item(113, 14).Value = -26213.166  # N113: -25774.143 -> -26213.166
$ws.Cells.Item(122, 8).Value = 3760.125  # H122: 4010.6086 -> 3760.125
$ws.Cells.Item(122, 9).Value = 2702.4211  # I122: 2963.7222 -> 2702.4211
$ws.Cells.Item(122, 11).Value = 8107.263300000001  # K122: 8891.1666 -> 8107.263300000001
$ws.Cells.Item(122, 13).Value = -5657.263300000001  # M122: -6441.1666 -> -5657.263300000001
$ws.Cells.Item(132, 8).Value = 2698.4666  # H132: 2814.1538 -> 2698.4666
$ws.Cells.Item(132, 9).Value = 2690.16  # I132: 2848.2 -> 2690.16
$ws.Cells.Item(132, 10).Value = 2740  # J132: 2700.6667 -> 2740
$ws.Cells.Item(132, 11).Value = 8070.48  # K132: 8544.599999999999 -> 8070.48
$ws.Cells.Item(132, 12).Value = 8220  # L132: 8102.000100000001 -> 8220
$ws.Cells.Item(132, 13).Value = -5540.48  # M132: -6014.599999999999 -> -5540.48
$ws.Cells.Item(132, 14).Value = -13280  # N132: -13162.0001 -> -13280

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1693.196  # H122: 1687.6735 -> 1693.196
$ws.Cells.Item(122, 9).Value = 1627.08  # I122: 1687.6735 -> 1627.08
$ws.Cells.Item(122, 10).Value = 4999  # J122: 0 -> 4999
$ws.Cells.Item(122, 11).Value = 4881.24  # K122: 5063.020500000001 -> 4881.24
$ws.Cells.Item(122, 12).Value = 14997  # L122: 0 -> 14997
$ws.Cells.Item(122, 13).Value = -2431.24  # M122: -2613.020500000001 -> -2431.24
$ws.Cells.Item(122, 14).Value = -19897  # N122: None -> -19897
$ws.Cells.Item(132, 8).Value = 6150.647  # H132: 7022.2856 -> 6150.647
$ws.Cells.Item(132, 9).Value = 1205  # I132: 1085.1428 -> 1205
$ws.Cells.Item(132, 10).Value = 11714.5  # J132: 12959.429 -> 11714.5
$ws.Cells.Item(132, 11).Value = 3615  # K132: 3255.4284 -> 3615
$ws.Cells.Item(132, 12).Value = 35143.5  # L132: 38878.287 -> 35143.5
$ws.Cells.Item(132, 13).Value = -1085  # M132: -725.4284000000002 -> -1085
$ws.Cells.Item(132, 14).Value = -40203.5  # N132: -43938.287 -> -40203.5
